$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$data = @(
    @(8, 13, 5, 7),
    @(3, 13, 5, 7),
    @(3, 16, 5, 4),
    @(5, 12, 4, 8),
    @(4, 12, 5, 8),
    @(8, 16, 5, 4),
    @(2, 14, 4, 6),
    @(3, 12, 2, 8),
    @(7, 17, 6, 3),
    @(3, 5, 7, 15),
    @(5, 5, 6, 15),
    @(5, 8, 6, 12),
    @(6, 5, 4, 15),
    @(5, 2, 2, 18),
    @(6, 17, 7, 3),
    @(3, 18, 2, 2),
    @(4, 2, 2, 18),
    @(7, 16, 6, 4),
    @(4, 5, 3, 15),
    @(2, 16, 5, 4),
    @(4, 7, 5, 13),
    @(5, 12, 4, 8),
    @(4, 18, 6, 2),
    @(4, 6, 3, 14),
    @(6, 6, 5, 14),
    @(2, 16, 3, 4),
    @(6, 6, 8, 14),
    @(4, 1, 5, 19),
    @(2, 14, 3, 6),
    @(4, 7, 3, 13),
    @(4, 6, 2, 14),
    @(5, 6, 9, 14),
    @(3, 7, 1, 13)
)

$startRow = 1370
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

$lastRow = $startRow + $data.Count - 1
$topLeftCellRow = $lastRow
$selCell = "A" + ($lastRow + 1)

$excel.ActiveWindow.ScrollRow = $topLeftCellRow
$ws.Range($selCell).Select()
